# Update betclever_predictions workbook with the latest predictions.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Home win"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")

$rows = @(
    @("26-12-2024 15:00", "ENGLAND", "CHAMPIONSHIP", "Watford - Portsmouth", 70, 1.85),
    @("26-12-2024 15:00", "ENGLAND", "LEAGUE ONE", "Huddersfield - Stockport County", 70, 2.1),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Dagenham & Redbridge - Braintree", 73.3, 1.8),
    @("26-12-2024 17:30", "ENGLAND", "NATIONAL LEAGUE", "Oldham - York", 80, 2.2),
    @("26-12-2024 13:00", "ALGERIA", "LIGUE 2", "GC Mascara - Oued Sly", 80, 1.7),
    @("26-12-2024 13:00", "ENGLAND", "NON LEAGUE PREMIER - NORTHERN", "Stockton Town - Hebburn Town", 80, 1.77)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "Draw"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")

$rows = @(
    @("25-12-2024 13:00", "TUNISIA", "LIGUE 2", "Jerba - EO Sidi Bouzid", 66.7, 2.65),
    @("25-12-2024 17:30", "WORLD", "GULF CUP OF NATIONS", "Bahrain - Iraq", 60, 2.9),
    @("26-12-2024 15:00", "ENGLAND", "LEAGUE TWO", "Bromley - Newport County", 60, 4.1),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Tamworth - Forest Green", 60, 3.6),
    @("26-12-2024 14:00", "ITALY", "SERIE B", "Reggiana - Juve Stabia", 60, 3),
    @("26-12-2024 13:00", "ALGERIA", "LIGUE 2", "Témouchent - WA Mostaganem", 60, 3),
    @("26-12-2024 15:00", "EGYPT", "PREMIER LEAGUE", "National Bank Of Egypt - Ceramica Cleopatra", 66.7, 2.9),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Eastbourne Borough - Salisbury", 60, 3.5),
    @("26-12-2024 13:00", "ENGLAND", "NON LEAGUE PREMIER - SOUTHERN SOUTH", "Sholing - Wimborne Town", 60, 3.3)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "Btts"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")

$rows = @(
    @("26-12-2024 17:30", "ENGLAND", "CHAMPIONSHIP", "Derby - West Brom", 75, 1.83),
    @("26-12-2024 15:00", "ENGLAND", "CHAMPIONSHIP", "Oxford United - Cardiff", 76, 1.73),
    @("26-12-2024 15:00", "ENGLAND", "CHAMPIONSHIP", "Preston - Hull City", 76.7, 1.83),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Aldershot Town - Woking", 76.7, 1.75),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Gateshead - Hartlepool", 76.7, 1.7),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Tamworth - Forest Green", 100, 1.8),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Enfield Town - Boreham Wood", 90, 2),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Weymouth - Torquay", 80, 1.75),
    @("26-12-2024 15:00", "NORTHERN-IRELAND", "PREMIERSHIP", "Crusaders FC - Cliftonville FC", 76.7, 1.7),
    @("26-12-2024 14:30", "WALES", "PREMIER LEAGUE", "Aberystwyth Town - Bala Town", 76.7, 1.85)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "Over_Under"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")

$rows = @(
    @("26-12-2024 15:00", "SCOTLAND", "PREMIERSHIP", "Dundee - Ross County", 85, 1.8, 60, 3),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Gateshead - Hartlepool", 75, 1.67, 68.8, 2.62),
    @("26-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Radcliffe - Chorley", 86.7, 1.7, 53.3, 2.88),
    @("26-12-2024 15:00", "ENGLAND", "NON LEAGUE PREMIER - SOUTHERN CENTRAL", "Kettering Town - Biggleswade Town", 93.3, 1.7, 53.3, $null),
    @("26-12-2024 15:00", "NORTHERN-IRELAND", "PREMIERSHIP", "Crusaders FC - Cliftonville FC", 80, 1.91, 50, 3.3)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -eq $null) {
        $ws.Cells.Item($r, 8).Value = ""
    } else {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "Away Win"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Away Win")

$rows = @(
    @("26-12-2024 13:00", "ENGLAND", "NON LEAGUE PREMIER - ISTHMIAN", "Bowers & Pitsea - Hashtag United", 73.3, 2.05),
    @("26-12-2024 15:00", "ENGLAND", "NON LEAGUE PREMIER - ISTHMIAN", "Hastings United - Dover", 70, 1.7)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
